$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)

# Update the title text
$title = $s.Shapes.Item("Title 1")
$title.TextFrame.TextRange.Text = "Conclusion / Lessons-Learned"

# Populate the body content placeholder with conclusions / lessons-learned bullets
$content = $s.Shapes.Item("Content Placeholder 2")
$tr = $content.TextFrame.TextRange
$tr.Text = "Cloud integration was challenging, requiring lots of troubleshooting/problem isolation."
$tr.InsertAfter("`rLesson-learned: start early as this task takes several days")
$tr.InsertAfter("`rLesson-learned: don" + [char]0x2019 + "t hard-code IP" + [char]0x2019 + "s in the scripts or config files")
$tr.InsertAfter("`rSome tasks may take collaboration across the team to resolve")
$tr.InsertAfter("`rLesson-learned: lean on your teammates for their expertise and insights; escalate issues before they become problems")

$tr.Paragraphs(1,1).LanguageID = "en-US"

$para2 = $tr.Paragraphs(2,1)
$para2.LanguageID = "en-US"
$para2.IndentLevel = 2

$para3 = $tr.Paragraphs(3,1)
$para3.LanguageID = "en-US"
$para3.IndentLevel = 2

$para4 = $tr.Paragraphs(4,1)
$para4.LanguageID = "en-US"

$para5 = $tr.Paragraphs(5,1)
$para5.LanguageID = "en-US"
$para5.IndentLevel = 2
